$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 22/23: Uniswap <-> Polygon swap (name + link) ---
$ws.Range("B22").Value = "Polygon"
$ws.Range("B23").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("C23").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"

# --- Column D (Price) -------------------------------------------------
# Excel auto-coerces numeric-looking strings (e.g. "1.00", "0.585") into
# numbers on assignment. The source data stores these as literal text, so
# force each contiguous block to Text format before writing the values,
# then restore the default "Normal" style afterwards (matching the
# original, un-styled cells).
$ws.Range("D2:D8").NumberFormat = "@"
$ws.Range("D11:D13").NumberFormat = "@"
$ws.Range("D15:D19").NumberFormat = "@"
$ws.Range("D21:D26").NumberFormat = "@"
$ws.Range("D28:D32").NumberFormat = "@"
$ws.Range("D34:D38").NumberFormat = "@"
$ws.Range("D42:D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D2").Value = "70.567.16"
$ws.Range("D3").Value = "3.535.73"
$ws.Range("D4").Value = "1.00"
$ws.Range("D5").Value = "607.47"
$ws.Range("D6").Value = "173.67"
$ws.Range("D7").Value = "0.617"
$ws.Range("D8").Value = "3.532.00"
$ws.Range("D11").Value = "6.75"
$ws.Range("D12").Value = "0.585"
$ws.Range("D13").Value = "47.45"
$ws.Range("D15").Value = "4.103.78"
$ws.Range("D16").Value = "626.27"
$ws.Range("D17").Value = "8.42"
$ws.Range("D18").Value = "70.614.52"
$ws.Range("D19").Value = "3.541.91"
$ws.Range("D21").Value = "17.42"
$ws.Range("D22").Value = "0.886"
$ws.Range("D23").Value = "9.93"
$ws.Range("D24").Value = "15.92"
$ws.Range("D25").Value = "96.94"
$ws.Range("D26").Value = "3.85"
$ws.Range("D28").Value = "2.61"
$ws.Range("D29").Value = "9.23"
$ws.Range("D30").Value = "33.44"
$ws.Range("D31").Value = "8.46"
$ws.Range("D32").Value = "3.10"
$ws.Range("D34").Value = "7.03"
$ws.Range("D35").Value = "567.87"
$ws.Range("D36").Value = "3.66"
$ws.Range("D37").Value = "10.77"
$ws.Range("D38").Value = "57.47"
$ws.Range("D42").Value = "0.0454"
$ws.Range("D43").Value = "0.328"
$ws.Range("D44").Value = "3.330.56"
$ws.Range("D45").Value = "3.04"
$ws.Range("D46").Value = "0.0₃0715"
$ws.Range("D47").Value = "33.08"
$ws.Range("D48").Value = "2.66"
$ws.Range("D50").Value = "133.60"
$ws.Range("D2:D8").Style = "Normal"
$ws.Range("D11:D13").Style = "Normal"
$ws.Range("D15:D19").Style = "Normal"
$ws.Range("D21:D26").Style = "Normal"
$ws.Range("D28:D32").Style = "Normal"
$ws.Range("D34:D38").Style = "Normal"
$ws.Range("D42:D48").Style = "Normal"
$ws.Range("D50").Style = "Normal"

# --- Column E (Volume 1h %) --------------------------------------------
# These are already non-numeric text (leading/trailing spaces, "%"), so a
# plain assignment keeps them as text exactly as in the source.
$ws.Range("E2").Value = "  +2.07%  "
$ws.Range("E3").Value = "  +0.98%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("E5").Value = "  +4.49%  "
$ws.Range("E6").Value = "  +0.39%  "
$ws.Range("E7").Value = "  -0.20%  "
$ws.Range("E8").Value = "  +1.04%  "
$ws.Range("E9").Value = "  -0.08%  "
$ws.Range("E10").Value = "  +6.62%  "
$ws.Range("E11").Value = "  +0.74%  "
$ws.Range("E12").Value = "  -1.46%  "
$ws.Range("E13").Value = "  +1.59%  "
$ws.Range("E14").Value = "  +1.98%  "
$ws.Range("E15").Value = "  +0.99%  "
$ws.Range("E16").Value = "  -7.07%  "
$ws.Range("E17").Value = "  -2.96%  "
$ws.Range("E18").Value = "  +2.13%  "
$ws.Range("E19").Value = "  +1.14%  "
$ws.Range("E20").Value = "  -1.69%  "
$ws.Range("E21").Value = "  +0.21%  "
$ws.Range("E22").Value = "  -1.41%  "
$ws.Range("E23").Value = "  -10.93%  "
$ws.Range("E24").Value = "  -1.07%  "
$ws.Range("E25").Value = "  -0.55%  "
$ws.Range("E26").Value = "  -0.15%  "
$ws.Range("E27").Value = "  +0.10%  "
$ws.Range("E28").Value = "  -1.39%  "
$ws.Range("E29").Value = "  -1.66%  "
$ws.Range("E30").Value = "  +1.77%  "
$ws.Range("E31").Value = "  -2.46%  "
$ws.Range("E32").Value = "  -2.50%  "
$ws.Range("E33").Value = "  -1.37%  "
$ws.Range("E34").Value = "  -2.99%  "
$ws.Range("E35").Value = "  -5.06%  "
$ws.Range("E36").Value = "  +2.17%  "
$ws.Range("E37").Value = "  -0.65%  "
$ws.Range("E38").Value = "  +0.74%  "
$ws.Range("E39").Value = "  -1.80%  "
$ws.Range("E40").Value = "  +0.15%  "
$ws.Range("E41").Value = "  +5.74%  "
$ws.Range("E42").Value = "  +4.11%  "
$ws.Range("E43").Value = "  -1.72%  "
$ws.Range("E44").Value = "  -2.46%  "
$ws.Range("E45").Value = "  +5.46%  "
$ws.Range("E46").Value = "  +1.31%  "
$ws.Range("E47").Value = "  -0.55%  "
$ws.Range("E48").Value = "  +2.49%  "
$ws.Range("E49").Value = "  -2.41%  "
$ws.Range("E50").Value = "  +0.33%  "
$ws.Range("E51").Value = "  -0.21%  "
